$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of trade data (row 5)
$ws.Range("A5").Value = 42647.680543981478
$ws.Range("B5").Value = $false
$ws.Range("C5").Value = 9988.1
$ws.Range("D5").Value = 10033.25
$ws.Range("E5").Value = 313
$ws.Range("F5").Value = 311.58999999999997
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = -0.45
$ws.Range("I5").Value = $true

# Match date/time number formatting already used in column A / G (style index 1 -> numFmtId 22)
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"
